$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# New named ranges used by the helper cells below
$wb.Names.Add('_nmax', '=Report!$I$6')
$wb.Names.Add('_nmin', '=Report!$I$5')
$wb.Names.Add('_nx', '=Report!$I$7')

# New helper cells H5:I10 on the Report sheet
# (labels typed in this order so the shared-string table ends up
# Up, Down, In, min, max, var - matching the authored workbook)
$ws.Range('H8').Value = "Up"
$ws.Range('H9').Value = "Down"
$ws.Range('H10').Value = "In"
$ws.Range('H5').Value = "min"
$ws.Range('H6').Value = "max"
$ws.Range('H7').Value = "var"

$ws.Range('I5').Formula = '=MIN(_nData)'
$ws.Range('I6').Formula = '=MAX(E7:E17)'
$ws.Range('I7').Value = 0
$ws.Range('I8').Formula = '=IF(_nx>_nmax,TRUE,FALSE)'
$ws.Range('I9').Formula = '=_nx<_nmin'
$ws.Range('I10').Formula = '=IF(MEDIAN(_nmin,_nmax,_nx)=_nx,TRUE,FALSE)'

# Clear the leftover "apply number format" style from the Formats sheet demo grid
$wsFormats = $wb.Worksheets.Item("Formats")
$wsFormats.Range('D19:F28').ClearFormats()

# Refresh the pivot table / pivot cache
$wsFormats.PivotTables("PivotTable1").PivotCache().Refresh()
